$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" "41.620.91"
$ws.Range("E2").Value = "  -1.08%  "
Set-TextValue "D3" "2.230.63"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.24%  "
Set-TextValue "D5" "252.57"
$ws.Range("E5").Value = "  +8.40%  "
Set-TextValue "D6" "0.630"
$ws.Range("E6").Value = "  -0.74%  "
Set-TextValue "D7" "71.27"
$ws.Range("E7").Value = "  +1.69%  "
$ws.Range("E8").Value = "  -0.08%  "
Set-TextValue "D9" "0.572"
$ws.Range("E9").Value = "  +2.01%  "
Set-TextValue "D10" "42.74"
$ws.Range("E10").Value = "  +20.64%  "
$ws.Range("E11").Value = "  -2.69%  "
Set-TextValue "D12" "58.93"
$ws.Range("E12").Value = "  +1.08%  "
Set-TextValue "D13" "0.105"
$ws.Range("E13").Value = "  +0.01%  "
Set-TextValue "D14" "7.01"
$ws.Range("E14").Value = "  +2.84%  "
Set-TextValue "D15" "2.552.66"
$ws.Range("E15").Value = "  -0.87%  "
Set-TextValue "D16" "14.98"
$ws.Range("E16").Value = "  -0.35%  "
Set-TextValue "D17" "0.853"
$ws.Range("E17").Value = "  -0.99%  "
Set-TextValue "D18" "2.224.02"
$ws.Range("E18").Value = "  -0.87%  "
Set-TextValue "D19" "41.596.07"
$ws.Range("E19").Value = "  -1.01%  "
Set-TextValue "D20" "0.0₃0968"
$ws.Range("E20").Value = "  -1.35%  "
Set-TextValue "D21" "6.22"
$ws.Range("E21").Value = "  -0.40%  "
Set-TextValue "D22" "73.07"
$ws.Range("E22").Value = "  -0.33%  "
Set-TextValue "D23" "2.29"
$ws.Range("E23").Value = "  +11.95%  "
Set-TextValue "D24" "234.96"
$ws.Range("E24").Value = "  -0.96%  "
Set-TextValue "D25" "3.85"
$ws.Range("E25").Value = "  +5.75%  "
$ws.Range("E26").Value = "  +0.16%  "
Set-TextValue "D27" "2.51"
$ws.Range("E27").Value = "  +6.46%  "
Set-TextValue "D28" "10.36"
$ws.Range("E28").Value = "  +3.12%  "
$ws.Range("E29").Value = "  +1.83%  "
Set-TextValue "D30" "171.32"
$ws.Range("E30").Value = "  +2.01%  "
Set-TextValue "D31" "20.67"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D33" "0.125"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D34" "5.60"
$ws.Range("E34").Value = "  +3.26%  "
$ws.Range("E35").Value = "  +0.70%  "
Set-TextValue "D36" "26.68"
$ws.Range("E36").Value = "  +20.38%  "
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("E38").Value = "  +11.28%  "
Set-TextValue "D39" "0.0288"
$ws.Range("E39").Value = "  +7.81%  "
$ws.Range("E40").Value = "  +2.24%  "
Set-TextValue "D41" "69.80"
$ws.Range("E41").Value = "  +2.86%  "
Set-TextValue "D42" "6.04"
$ws.Range("E42").Value = "  -0.23%  "
Set-TextValue "D43" "12.04"
$ws.Range("E43").Value = "  +19.33%  "
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("E45").Value = "  +10.47%  "
Set-TextValue "D46" "8.85"
$ws.Range("E46").Value = "  -2.87%  "
Set-TextValue "D47" "4.79"
$ws.Range("E47").Value = "  +9.45%  "
Set-TextValue "D48" "0.102"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("E50").Value = "  +6.97%  "
Set-TextValue "D51" "1.20"
$ws.Range("E51").Value = "  +1.70%  "
